$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 2023-07-19 / 15011 row currently sitting in row 2 moves down to row 4.
# Use Copy (not Value=) so the date-like text stays text instead of being
# reinterpreted as a date serial number.
$ws.Range("A2:B2").Copy($ws.Range("A4:B4"))

# Remove the old row 1 (2023-07-18 / 13366) and the now-duplicated row 2.
$ws.Range("A1:B2").ClearContents()
